$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F6").Value = 9161
$wsExhibit.Range("F10").Value = 1058
$wsExhibit.Range("F12").Value = 53

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F8").Value = 9161
$wsAll.Range("F12").Value = 1058
$wsAll.Range("F14").Value = 53
